# Generate Report for Archive
#
# 1. Flip the localization status text from "Ready for handoff" to
#    "In Translation" everywhere it appears:
#      - Overview sheet: E2:F3 (the zh-cn / de-de status columns)
#      - zh-cn sheet:     C2:C3 (Status column)
#      - de-de sheet:     C2:C3 (Status column)
#
# 2. Narrow the "Status" columns (Overview!E:F, zh-cn!C, de-de!C) from
#    ~17.22 characters down to ~13.41 characters.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Update the status text -------------------------------------------------
$wsOverview.Range("E2:F3").Value = "In Translation"
$wsZhCn.Range("C2:C3").Value     = "In Translation"
$wsDeDe.Range("C2:C3").Value     = "In Translation"

# --- Resize the status columns ----------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth     = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth     = 12.5
